$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 86; this shifts existing rows 86-208 down to 87-209
# (and carries D-column's date style down with them automatically).
$ws.Rows.Item(86).EntireRow.Insert()

# Populate the newly inserted row 86 with the new data record.
$ws.Range("A86").Value = 8
$ws.Range("B86").Value = "Terminal La Palmera de La Serena"
$ws.Range("C86").Value = "Coquimbo"
$ws.Range("D86").Value = 44994
$ws.Range("E86").Value = 4
$ws.Range("F86").Value = 100112044
$ws.Range("G86").Value = "Perejil"
$ws.Range("H86").Value = "Sin especificar"
$ws.Range("I86").Value = "Primera"
$ws.Range("J86").Value = 2000
$ws.Range("K86").Value = 2000
$ws.Range("L86").Value = 2500
$ws.Range("M86").Value = 2250
$ws.Range("N86").Value = "`$/atado 1 a 1,5 kilos"
$ws.Range("O86").Value = "Provincia del Elquí"
$ws.Range("P86").Value = 1500
$ws.Range("Q86").Value = 1.5
$ws.Range("R86").Value = "Hortaliza"
